$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "展览" (Exhibitions) — F column "want-to-go" counters bump up, plus
# two cover-image URL refreshes on rows 31/32.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item('展览')
$ws1.Range('F2').Value = 203
$ws1.Range('F3').Value = 2478
$ws1.Range('F5').Value = 1842
$ws1.Range('F6').Value = 115
$ws1.Range('F7').Value = 329
$ws1.Range('F8').Value = 631
$ws1.Range('F9').Value = 3622
$ws1.Range('F10').Value = 1228
$ws1.Range('F11').Value = 1592
$ws1.Range('F12').Value = 32
$ws1.Range('F15').Value = 1448
$ws1.Range('F16').Value = 2
$ws1.Range('F17').Value = 1807
$ws1.Range('F20').Value = 13
$ws1.Range('F21').Value = 480
$ws1.Range('F22').Value = 1563
$ws1.Range('F28').Value = 281
$ws1.Range('F30').Value = 4392
$ws1.Range('F31').Value = 60
$ws1.Range('I31').Value = '//i2.hdslb.com/bfs/openplatform/202410/CXUc87f81729246062802.jpeg'
$ws1.Range('F32').Value = 60
$ws1.Range('I32').Value = '//i2.hdslb.com/bfs/openplatform/202410/CXUc87f81729246062802.jpeg'
$ws1.Range('F33').Value = 12
$ws1.Range('F35').Value = 156
$ws1.Range('F37').Value = 1226
$ws1.Range('F38').Value = 957

# ---------------------------------------------------------------------------
# Sheet "演出" (Performances) — row 2 is a full event swap (new concert),
# plus assorted F-column counter bumps.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item('演出')
$ws2.Range('C2').Value = '上海·ROOKiEZ is PUNK`D 「Reignite Youth （重燃青春）」2024 CHINA Tour '
$ws2.Range('D2').Value = '虹许路731号4号楼 THE BOXX•城市乐园'
$ws2.Range('E2').Value = '2024.10.18 20:30-10.18 22:00'
$ws2.Range('F2').Value = 81
$ws2.Range('H2').Value = 'https://show.bilibili.com/platform/detail.html?id=91376'
$ws2.Range('I2').Value = '//i1.hdslb.com/bfs/openplatform/202408/pZdI02BJ1724735899119.jpeg'

$ws2.Range('F5').Value = 30
$ws2.Range('F22').Value = 156
$ws2.Range('F24').Value = 196
$ws2.Range('F28').Value = 4
$ws2.Range('F29').Value = 69
$ws2.Range('F35').Value = 454
$ws2.Range('F39').Value = 4
$ws2.Range('F40').Value = 34
$ws2.Range('F43').Value = 92
$ws2.Range('F47').Value = 33
$ws2.Range('F48').Value = 33

# ---------------------------------------------------------------------------
# Sheet "本地生活" (Local Life) — F-column counter bumps, then the two
# "NIJISANJI EN X KAKACODE" rows (old rows 14 & 15) are removed outright;
# what was row 16 (WIND BREAKER x animate cafe) shifts up to become row 14
# and its want-to-go / price numbers are refreshed.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item('本地生活')
$ws3.Range('F3').Value = 2551
$ws3.Range('F4').Value = 2560
$ws3.Range('F5').Value = 9588
$ws3.Range('F6').Value = 166
$ws3.Range('F9').Value = 400
$ws3.Range('F10').Value = 3035
$ws3.Range('F11').Value = 543
$ws3.Range('F12').Value = 832

$ws3.Range('A14:A15').EntireRow.Delete()

$ws3.Range('A14').Value = 13
$ws3.Range('F14').Value = 280
$ws3.Range('G14').Value = 30

# ---------------------------------------------------------------------------
# Sheet "全部类型" (All Types) — same counter bumps as above (this sheet is
# an independently-maintained composite view), plus a cascading content
# shift across rows 37-41 (each row's event slides into the next slot) and
# a couple of trailing cover/count refreshes.
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item('全部类型')
$ws4.Range('F2').Value = 2551
$ws4.Range('F3').Value = 166
$ws4.Range('F4').Value = 203
$ws4.Range('F5').Value = 2478
$ws4.Range('F7').Value = 543
$ws4.Range('F8').Value = 832
$ws4.Range('F9').Value = 115
$ws4.Range('F10').Value = 329
$ws4.Range('F11').Value = 631
$ws4.Range('F12').Value = 3622
$ws4.Range('F13').Value = 1228
$ws4.Range('F17').Value = 1448
$ws4.Range('F22').Value = 1807
$ws4.Range('F24').Value = 13
$ws4.Range('F25').Value = 480
$ws4.Range('F27').Value = 1563
$ws4.Range('F29').Value = 156
$ws4.Range('F30').Value = 156
$ws4.Range('F32').Value = 196
$ws4.Range('F36').Value = 281

$ws4.Range('B37').Value = '2024-11-15'
$ws4.Range('C37').Value = '上海·“法国姐姐”乔伊丝·乔纳森《小意思》'
$ws4.Range('D37').Value = '高青西路777号 上海前滩31演艺中心'
$ws4.Range('E37').Value = '2024.11.15 19:30-11.15 21:00'
$ws4.Range('F37').Value = 5
$ws4.Range('G37').Value = 280
$ws4.Range('H37').Value = 'https://show.bilibili.com/platform/detail.html?id=91619'
$ws4.Range('I37').Value = '//i1.hdslb.com/bfs/openplatform/202408/VnZEk71H1725014748758.jpeg'

$ws4.Range('C38').Value = '上海·「WIND BREAKER × animate cafe」'
$ws4.Range('D38').Value = '西藏北路198号大悦城北座8楼N809-1 animate cafe上海店'
$ws4.Range('E38').Value = '2024.11.15 00:00-12.15 23:59'
$ws4.Range('F38').Value = 280
$ws4.Range('G38').Value = 30
$ws4.Range('H38').Value = 'https://show.bilibili.com/platform/detail.html?id=93422'
$ws4.Range('I38').Value = '//i0.hdslb.com/bfs/openplatform/202410/TGPx1EZW1728892799830.jpeg'

$ws4.Range('B39').Value = '2024-11-16'
$ws4.Range('C39').Value = '上海·变形金刚音乐会40周年特别版'
$ws4.Range('D39').Value = '高青西路777号 上海前滩31演艺中心'
$ws4.Range('E39').Value = '2024.11.16 19:30-11.16 21:30'
$ws4.Range('F39').Value = 63
$ws4.Range('G39').Value = 266
$ws4.Range('H39').Value = 'https://show.bilibili.com/platform/detail.html?id=90031'
$ws4.Range('I39').Value = '//i1.hdslb.com/bfs/openplatform/202409/5zTUqO9f1727061199503.jpeg'

$ws4.Range('C40').Value = '上海·趣元界·第三届ICG动漫游戏博览会'
$ws4.Range('D40').Value = '西藏南路1号 上海大世界'
$ws4.Range('E40').Value = '2024.11.16 10:00-11.17 17:00'
$ws4.Range('F40').Value = 4392
$ws4.Range('G40').Value = 59
$ws4.Range('H40').Value = 'https://show.bilibili.com/platform/detail.html?id=92846'
$ws4.Range('I40').Value = '//i1.hdslb.com/bfs/openplatform/202410/C1h14i9R1728540930986.jpeg'

$ws4.Range('B41').Value = '2024-11-22'
$ws4.Range('C41').Value = '上海·第五届长三角文博会上海国际插画艺术节'
$ws4.Range('D41').Value = '崧泽大道333号 国家会展中心'
$ws4.Range('E41').Value = '2024.11.22 09:00-11.25 15:00'
$ws4.Range('F41').Value = 60
$ws4.Range('G41').Value = 62
$ws4.Range('H41').Value = 'https://show.bilibili.com/platform/detail.html?id=92813'
$ws4.Range('I41').Value = '//i2.hdslb.com/bfs/openplatform/202410/CXUc87f81729246062802.jpeg'

$ws4.Range('F42').Value = 60
$ws4.Range('I42').Value = '//i2.hdslb.com/bfs/openplatform/202410/CXUc87f81729246062802.jpeg'
$ws4.Range('F43').Value = 454
$ws4.Range('F44').Value = 12
$ws4.Range('F46').Value = 156
$ws4.Range('F47').Value = 92
$ws4.Range('F50').Value = 33
$ws4.Range('F51').Value = 1226
$ws4.Range('F52').Value = 957
